# Updated cryptos list on Sun Jun  4 22:25:18 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row($row, $d, $e) {
    if ($d -ne $null) {
        # Force text format so numeric-looking strings (with leading/trailing
        # zeros, multiple dots, etc.) are preserved exactly like the source data.
        $ws.Range("D$row").NumberFormat = "@"
        $ws.Range("D$row").Value = $d
    }
    if ($e -ne $null) { $ws.Range("E$row").Value = $e }
}

Set-Row 2  "27.312.71"    "  +0.71%  "
Set-Row 3  "1.905.81"     $null
Set-Row 4  "0.9986"       "  -0.41%  "
Set-Row 5  "306.67"       "  -0.09%  "
Set-Row 6  "0.9980"       "  -0.45%  "
Set-Row 7  "0.5419"       "  +3.99%  "
Set-Row 8  "0.3811"       "  +1.30%  "
Set-Row 9  "0.07310"      "  +0.54%  "
Set-Row 10 $null           "  +4.27%  "
Set-Row 11 "0.9044"       "  +0.57%  "
Set-Row 12 "0.08191"      "  -0.07%  "
Set-Row 13 "95.70"        "  -0.56%  "
Set-Row 14 "5.358"        "  +1.12%  "
Set-Row 15 "0.9982"       "  -0.47%  "
Set-Row 16 "14.87"        "  +2.04%  "
Set-Row 17 "0.000008661"  "  +0.69%  "
Set-Row 18 "1.357.92"     "  -28.62%  "
Set-Row 19 "0.9976"       "  -0.58%  "
Set-Row 20 "27.212.29"    "  +0.26%  "
Set-Row 21 "5.055"        "  -0.64%  "
Set-Row 22 $null           "  +1.05%  "
Set-Row 23 "6.514"        "  +1.58%  "

# Row 24 and 25 swap content: Monero <-> LidoDAOToken, then values updated
$ws.Range("B24").Value = "LidoDAOToken"
$ws.Range("C24").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
Set-Row 24 "2.314" "  +0.25%  "

$ws.Range("B25").Value = "Monero"
$ws.Range("C25").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-Row 25 "148.50" "  +0.00%  "

Set-Row 26 "18.36"   "  +0.92%  "
Set-Row 27 "1.755"   "  +0.80%  "
Set-Row 28 "116.82"  "  +1.44%  "
Set-Row 29 "4.850"   "  +1.20%  "
Set-Row 30 "4.680"   "  -3.86%  "
Set-Row 31 "0.09191" "  -0.15%  "
Set-Row 32 "0.8304"  "  +4.80%  "
Set-Row 33 "0.05078" "  +0.91%  "
Set-Row 34 "1.228"   "  +1.01%  "
Set-Row 35 "3.011"   "  +1.20%  "
Set-Row 36 "3.326"   "  -2.95%  "
Set-Row 37 "2.711"   "  +3.93%  "
Set-Row 38 "0.6023"  "  +5.33%  "
Set-Row 39 $null      "  +0.43%  "
Set-Row 40 "1.077"   "  +0.27%  "
Set-Row 41 "9.282"   "  +2.79%  "
Set-Row 42 "6.669"   "  +1.75%  "

# Row 43 and 44 swap content: Decentraland <-> Quant, then values updated
$ws.Range("B43").Value = "Quant"
$ws.Range("C43").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
Set-Row 43 "116.12" "  -0.35%  "

$ws.Range("B44").Value = "Decentraland"
$ws.Range("C44").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
Set-Row 44 "0.5182" "  +6.73%  "

Set-Row 45 "0.1533"  "  +1.27%  "
Set-Row 46 "10.23"   "  +1.86%  "
Set-Row 47 "0.9972"  "  -0.56%  "
Set-Row 48 $null      "  +1.28%  "
Set-Row 49 "38.11"   "  -0.28%  "
Set-Row 50 "0.06089" "  +2.67%  "
Set-Row 51 "63.56"   "  -0.05%  "
